$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    ,("ECs","Hsp90aa1","Fgfr3","ECs",3,1,63.73255033333334,191.197651,0.09718402276460011,0.1059076069828809,3,1,2.330840333333333,6.992521,0.6715345129768794,0.7003397275969581,148.5503988631301,1336.953589768171,0.0652624253963597,0.07417130462483654)
    ,("ECs","Hsp90aa1","Fgfr3","FAPs",3,1,63.73255033333334,191.197651,0.09718402276460011,0.1059076069828809,3,1,0.657666,1.972998,0.1894790521235985,0.1976066831789769,41.914731447522,377.232583027698,0.01841433651499465,0.02092805093930976)
    ,("ECs","Hsp90aa1","Fgfr3","MuSCs",3,1,63.73255033333334,191.197651,0.09718402276460011,0.1059076069828809,2,1,0.42828,0.85656,0.123391035029171,0.08578923067523865,27.29537665676,163.77225994056,0.01199163715722252,0.00908573212571689)
    ,("ECs","Hsp90aa1","Fgfr3","Resolving-Mac",3,1,63.73255033333334,191.197651,0.09718402276460011,0.1059076069828809,1,0.3333333333333333,0.05413033333333334,0.162391,0.01559539987035126,0.01626435854882633,3.449864193726778,31.048777743541,0.001515623696023258,0.001722519293017759)
    ,("FAPs","Hsp90aa1","Fgfr3","ECs",3,1,125.4169656666667,376.250897,0.1912449004891238,0.2084117242969288,3,1,2.330840333333333,6.992521,0.6715345129768794,0.7003397275969581,292.3269220601485,2630.942298541337,0.1284275511092755,0.1459590102221235)
    ,("FAPs","Hsp90aa1","Fgfr3","FAPs",3,1,125.4169656666667,376.250897,0.1912449004891238,0.2084117242969288,3,1,0.657666,1.972998,0.1894790521235985,0.1976066831789769,82.482474142134,742.342267279206,0.03623690246815109,0.04118354957392751)
    ,("FAPs","Hsp90aa1","Fgfr3","MuSCs",3,1,125.4169656666667,376.250897,0.1912449004891238,0.2084117242969288,2,1,0.42828,0.85656,0.123391035029171,0.08578923067523865,53.71357805572,322.28146833432,0.02359790621540379,0.01787948149113347)
    ,("FAPs","Hsp90aa1","Fgfr3","Resolving-Mac",3,1,125.4169656666667,376.250897,0.1912449004891238,0.2084117242969288,1,0.3333333333333333,0.05413033333333334,0.162391,0.01559539987035126,0.01626435854882633,6.78886215719189,61.099759414727,0.00298254069629342,0.00338968300974439)
    ,("Inflammatory-Mac","Hsp90aa1","Fgfr3","ECs",3,1,184.1540323333334,552.4620970000001,0.2808114468489298,0.3060180830465028,3,1,2.330840333333333,6.992521,0.6715345129768794,0.7003397275969581,429.2336461085042,3863.102814976538,0.1885745781980289,0.2143166209205311)
    ,("Inflammatory-Mac","Hsp90aa1","Fgfr3","FAPs",3,1,184.1540323333334,552.4620970000001,0.2808114468489298,0.3060180830465028,3,1,0.657666,1.972998,0.1894790521235985,0.1976066831789769,121.111845828534,1090.006612456806,0.05320788677439148,0.06047121838360814)
    ,("Inflammatory-Mac","Hsp90aa1","Fgfr3","MuSCs",3,1,184.1540323333334,552.4620970000001,0.2808114468489298,0.3060180830465028,2,1,0.42828,0.85656,0.123391035029171,0.08578923067523865,78.86948896772002,473.2169338063201,0.03464961507472848,0.02625305591727077)
    ,("Inflammatory-Mac","Hsp90aa1","Fgfr3","Resolving-Mac",3,1,184.1540323333334,552.4620970000001,0.2808114468489298,0.3060180830465028,1,0.3333333333333333,0.05413033333333334,0.162391,0.01559539987035126,0.01626435854882633,9.96831915488078,89.71487239392702,0.004379366801780948,0.004977187825092833)
    ,("MuSCs","Hsp90aa1","Fgfr3","ECs",2,1,162.052406,324.104812,0.2471092813859239,0.1795271274047008,3,1,2.330840333333333,6.992521,0.6715345129768794,0.7003397275969581,377.7182840185087,2266.309704111052,0.165942410927563,0.1257299795028725)
    ,("MuSCs","Hsp90aa1","Fgfr3","FAPs",2,1,162.052406,324.104812,0.2471092813859239,0.1795271274047008,3,1,0.657666,1.972998,0.1894790521235985,0.1976066831789769,106.576357644396,639.4581458663761,0.04682203240794844,0.03547576018709254)
    ,("MuSCs","Hsp90aa1","Fgfr3","MuSCs",2,1,162.052406,324.104812,0.2471092813859239,0.1795271274047008,2,1,0.42828,0.85656,0.123391035029171,0.08578923067523865,69.40380444168001,277.61521776672,0.0304910699955238,0.01540149414538484)
    ,("MuSCs","Hsp90aa1","Fgfr3","Resolving-Mac",2,1,162.052406,324.104812,0.2471092813859239,0.1795271274047008,1,0.3333333333333333,0.05413033333333334,0.162391,0.01559539987035126,0.01626435854882633,8.771950754248667,52.63170452549201,0.003853768054888629,0.002919893569350879)
    ,("Resolving-Mac","Hsp90aa1","Fgfr3","ECs",3,1,120.4365156666667,361.309547,0.1836503485114226,0.2001354582689865,3,1,2.330840333333333,6.992521,0.6715345129768794,0.7003397275969581,280.7182883219986,2526.464594897987,0.1233275473456523,0.1401628123265944)
    ,("Resolving-Mac","Hsp90aa1","Fgfr3","FAPs",3,1,120.4365156666667,361.309547,0.1836503485114226,0.2001354582689865,3,1,0.657666,1.972998,0.1894790521235985,0.1976066831789769,79.20700151243399,712.863013611906,0.03479789395811288,0.03954810409503898)
    ,("Resolving-Mac","Hsp90aa1","Fgfr3","MuSCs",3,1,120.4365156666667,361.309547,0.1836503485114226,0.2001354582689865,2,1,0.42828,0.85656,0.123391035029171,0.08578923067523865,51.58055092972,309.48330557832,0.02266080658629241,0.01716946699573269)
    ,("Resolving-Mac","Hsp90aa1","Fgfr3","Resolving-Mac",3,1,120.4365156666667,361.309547,0.1836503485114226,0.2001354582689865,1,0.3333333333333333,0.05413033333333334,0.162391,0.01559539987035126,0.01626435854882633,6.519268738541889,58.673418646877,0.002864100621365003,0.003255074851620466)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $r = $startRow + $i
    $vals = $rowData[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

Write-Output "Wrote $($rowData.Length) rows starting at row $startRow"
